$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>madrid</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$t.Cell(1,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>barcelona</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$t.Cell(2,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>l</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>azio</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$t.Cell(2,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>b</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>etis</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$t.Cell(3,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>b</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>ayern</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>')
$t.Cell(3,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>m</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>anchester</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>c</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>ity</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$t.Cell(4,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>porting</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>')
$t.Cell(4,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve">nter </w:t></w:r></w:p>')
$t.Cell(5,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>psg</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$t.Cell(5,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:kern w:val="0"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr><w:t>rsenal</w:t></w:r></w:p>')

Write-Output "done"